$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list: refreshed prices (col D) and 1h volume % (col E).
# Two pairs of rows (30/31 and 45/46) also swapped rank order, so their
# Coin/Link/Price/Volume (B/C/D/E) cells exchange contents below.
#
# All of these are plain-text cells in the source sheet (coinranking
# prices/percentages, not real numbers). We prefix the literal value with
# a quote-prefix (') so numeric-looking strings like "580.32" are not
# auto-converted to the Number type by the Value setter, then reset the
# cell Style back to Normal so no stray quote-prefix formatting remains
# (matches the un-styled inline-string cells in the original workbook).

$ws.Range('D2').Value = "'" + '62.206.92'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  -2.33%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '2.997.55'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  -2.61%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  -0.02%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '580.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  -1.38%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '146.69'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  -5.51%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  +0.00%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  -3.03%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'" + '2.993.87'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  -2.68%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  -5.29%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  -3.73%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  -2.34%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  -4.44%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'" + '34.56'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  -6.01%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  +1.51%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'" + '3.494.01'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  -2.56%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'" + '7.00'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  -2.25%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'" + '62.246.19'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  -2.15%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'" + '2.998.76'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  -2.54%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'" + '452.77'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  -3.48%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'" + '13.81'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  -3.25%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'" + '0.676'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  -3.97%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'" + '7.27'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  -2.73%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'" + '2.28'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  -6.14%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'" + '79.92'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  -0.59%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  -4.44%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'" + '10.03'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  -3.92%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  -0.06%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  +0.03%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('B30').Value = "'" + 'PancakeSwap'
$ws.Range('B30').Style = 'Normal'
$ws.Range('C30').Value = "'" + 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('C30').Style = 'Normal'
$ws.Range('D30').Value = "'" + '2.61'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  -2.05%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('B31').Value = "'" + 'NEARProtocol'
$ws.Range('B31').Style = 'Normal'
$ws.Range('C31').Value = "'" + 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('C31').Style = 'Normal'
$ws.Range('D31').Value = "'" + '7.11'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  -4.35%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  -2.18%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'" + '26.83'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  -0.88%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  -5.13%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  -2.25%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'" + '0.0₃0792'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  -4.25%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  -4.32%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'" + '2.11'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  -3.81%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  -0.61%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  -2.14%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  -9.98%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'" + '409.95'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  -5.23%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  -0.69%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'" + '0.275'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  -4.96%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('B45').Value = "'" + 'Maker'
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').Value = "'" + 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').Value = "'" + '2.758.10'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  -1.86%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').Value = "'" + 'VeChain'
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = "'" + 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = "'" + '0.0351'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  -2.60%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'" + '38.10'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  -4.20%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'" + '127.92'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  -1.69%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D50').Value = "'" + '0.107'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  -1.69%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'" + '23.64'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  -5.02%  '
$ws.Range('E51').Style = 'Normal'
